$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") on every data row (2-514) was refreshed by one day:
# 45179 (2023-09-10) -> 45180 (2023-09-11)
for ($r = 2; $r -le 514; $r++) {
    $ws.Cells.Item($r, 3).Value = 45180
}

# Row 5 (A 60950-2019) additionally gained a new observed species
# ("Siljansspindling"), which bumps several of its summary counters:
#   VU (K)          8  -> 9
#   Rödlistade (O) 19  -> 20
#   Hotade (P)      9  -> 10
#   Alla arter (Q) 20  -> 21
$ws.Range("K5").Value = 9
$ws.Range("O5").Value = 20
$ws.Range("P5").Value = 10
$ws.Range("Q5").Value = 21

# The species list in R5 gets "Siljansspindling" inserted right after
# "Rynkskinn" (keeping the same CRLF-separated layout as the original).
$ws.Range("R5").Value = "Violett guldvinge`r`nFjällfotad fingersvamp`r`nFjällfotad musseron`r`nGrangråticka`r`nLäderdoftande fingersvamp`r`nRynkskinn`r`nSiljansspindling`r`nSpricktaggsvamp`r`nSvartfjällig musseron`r`nTaggfingersvamp`r`nBlek fingersvamp`r`nDofttaggsvamp`r`nDruvfingersvamp`r`nFlattoppad klubbsvamp`r`nGultoppig fingersvamp`r`nOrange taggsvamp`r`nRosenticka`r`nRödbrun klubbdyna`r`nUllticka`r`nÄggvaxskivling`r`nGranriska"
